# Insert a new data row at row 426 (this pushes the existing rows 426-514
# down to 427-515, keeping rows 1-425 untouched) and populate the newly
# inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 426.
$ws.Rows.Item(426).Insert()

# Populate the new row 426 with the record's data.
$ws.Range("A426").Value = 5
$ws.Range("B426").Value = "Macroferia Regional de Talca"
$ws.Range("C426").Value = "Maule"
$ws.Range("D426").Value = 44641
$ws.Range("E426").Value = 7
$ws.Range("F426").Value = "Fruta"
$ws.Range("G426").Value = 100101
$ws.Range("H426").Value = "Berries"
$ws.Range("I426").Value = 100112025
$ws.Range("J426").Value = "Frutilla"
$ws.Range("K426").Value = "Sin especificar"
$ws.Range("L426").Value = "Primera"
$ws.Range("M426").Value = 200
$ws.Range("N426").Value = 7000
$ws.Range("O426").Value = 7000
$ws.Range("P426").Value = 7000
$ws.Range("Q426").Value = "`$/bandeja 7 kilos"
$ws.Range("R426").Value = "Región del Maule"
$ws.Range("S426").Value = 1000
$ws.Range("T426").Value = 7
